$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9934803247451782
$ws.Range("B1").Value = 2.02260422706604
$ws.Range("C1").Value = 5.33855676651001
$ws.Range("D1").Value = 1.067158937454224
$ws.Range("E1").Value = 0.858919620513916
